$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to stay text (the source workbook stores all prices as text)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.070.26"
$ws.Range("E2").Value = "  +1.93%  "

$ws.Range("D3").Value = "3.481.44"
$ws.Range("E3").Value = "  +2.19%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "416.33"
$ws.Range("E5").Value = "  +1.78%  "

$ws.Range("D6").Value = "131.63"

$ws.Range("D7").Value = "0.632"
$ws.Range("E7").Value = "  -0.77%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "0.734"
$ws.Range("E9").Value = "  -0.17%  "

$ws.Range("E10").Value = "  +8.43%  "

$ws.Range("D11").Value = "42.86"
$ws.Range("E11").Value = "  -1.16%  "

$ws.Range("D12").Value = "9.82"
$ws.Range("E12").Value = "  +5.07%  "

$ws.Range("D13").Value = "0.0000227"
$ws.Range("E13").Value = "  +2.38%  "

$ws.Range("D14").Value = "4.046.90"
$ws.Range("E14").Value = "  +2.41%  "

$ws.Range("E15").Value = "  -0.19%  "

$ws.Range("D16").Value = "20.67"
$ws.Range("E16").Value = "  -3.29%  "

$ws.Range("D17").Value = "3.478.21"
$ws.Range("E17").Value = "  +1.80%  "

$ws.Range("D18").Value = "12.72"
$ws.Range("E18").Value = "  +1.36%  "

$ws.Range("E19").Value = "  +0.46%  "

$ws.Range("D20").Value = "63.022.95"
$ws.Range("E20").Value = "  +1.81%  "

$ws.Range("D21").Value = "468.49"
$ws.Range("E21").Value = "  +4.29%  "

$ws.Range("D22").Value = "90.96"
$ws.Range("E22").Value = "  -0.75%  "

$ws.Range("E23").Value = "  +3.46%  "

$ws.Range("D24").Value = "13.31"

$ws.Range("D25").Value = "10.71"
$ws.Range("E25").Value = "  +14.83%  "

$ws.Range("D26").Value = "3.35"
$ws.Range("E26").Value = "  +1.60%  "

$ws.Range("D27").Value = "33.66"
$ws.Range("E27").Value = "  +1.44%  "

$ws.Range("D28").Value = "4.79"
$ws.Range("E28").Value = "  +0.09%  "

$ws.Range("D29").Value = "7.59"
$ws.Range("E29").Value = "  -0.75%  "

$ws.Range("D30").Value = "12.18"
$ws.Range("E30").Value = "  +1.24%  "

$ws.Range("E31").Value = "  -1.45%  "

$ws.Range("E32").Value = "  -0.57%  "

$ws.Range("E33").Value = "  -0.99%  "

$ws.Range("D34").Value = "41.17"
$ws.Range("E34").Value = "  -3.19%  "

$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("D36").Value = "58.48"
$ws.Range("E36").Value = "  +8.67%  "

$ws.Range("E37").Value = "  -2.42%  "

$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.12%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "3.08"
$ws.Range("E39").Value = "  +4.11%  "

$ws.Range("E40").Value = "  +7.37%  "

$ws.Range("E41").Value = "  -0.23%  "

$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "148.16"
$ws.Range("E42").Value = "  +3.05%  "

$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").Value = "4.48"
$ws.Range("E43").Value = "  +2.53%  "

$ws.Range("D44").Value = "0.322"
$ws.Range("E44").Value = "  +1.11%  "

$ws.Range("E45").Value = "  -1.26%  "

$ws.Range("E46").Value = "  +3.17%  "

$ws.Range("D47").Value = "0.0₃0585"
$ws.Range("E47").Value = "  +34.45%  "

$ws.Range("E48").Value = "  +11.79%  "

$ws.Range("D49").Value = "16.48"
$ws.Range("E49").Value = "  -1.08%  "

$ws.Range("D50").Value = "22.28"
$ws.Range("E50").Value = "  -1.37%  "

$ws.Range("D51").Value = "0.145"
$ws.Range("E51").Value = "  -3.26%  "
